$d = $word.ActiveDocument

# Locate the Subtitle paragraph ("Reflections on the implications ...")
# so we can insert the new Author paragraph directly after it.
$subtitleIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Reflections on the implications of desistance theory for long-term imprisonment*") {
        $subtitleIndex = $i
        break
    }
}

if ($subtitleIndex -eq 0) {
    throw "Could not locate the subtitle paragraph to anchor the new Author paragraph."
}

# Insert a new, empty paragraph right after the subtitle.
$subtitle = $d.Paragraphs.Item($subtitleIndex)
$subtitle.Range.InsertParagraphAfter()

# Populate the newly created paragraph with the author's name and apply
# the "Author" paragraph style.
$authorPara = $d.Paragraphs.Item($subtitleIndex + 1)
$authorPara.Range.Text = "Ben Jarman"
$authorPara.Style = "Author"
